$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Range('B114').Value = 3373723
$ws.Range('C114').Value = 'Uruguay Primera División'
$ws.Range('D114').Value = 'Uruguay Clausura'
$ws.Range('E114').Value = 44283.83333333334
$ws.Range('F114').Value = 'Cerro Largo'
$ws.Range('G114').Value = 'Defensor Sporting'
$ws.Range('H114').Value = 0
$ws.Range('I114').Value = 0
$ws.Range('J114').Value = 'D'
$ws.Range('K114').Value = 2.625
$ws.Range('L114').Value = 3.3
$ws.Range('M114').Value = 2.375
$ws.Range('N114').Value = 2.8
$ws.Range('O114').Value = 3
$ws.Range('P114').Value = 2.45
$ws.Range('Q114').Value = 0
$ws.Range('R114').Value = 2.125
$ws.Range('S114').Value = 1.75
$ws.Range('T114').Value = 2.25
$ws.Range('U114').Value = 1.925
$ws.Range('V114').Value = 1.925
$ws.Range('W114').Value = -1
$ws.Range('X114').Value = 2
$ws.Range('Y114').Value = -1
$ws.Range('Z114').Value = 0
$ws.Range('AA114').Value = -0
$ws.Range('AB114').Value = -1
$ws.Range('AC114').Value = 0.925

# Row 116
$ws.Range('B116').Value = 3373720
$ws.Range('C116').Value = 'Uruguay Primera División'
$ws.Range('D116').Value = 'Uruguay Clausura'
$ws.Range('E116').Value = 44283.83333333334
$ws.Range('F116').Value = 'Rentistas'
$ws.Range('G116').Value = 'Danubio'
$ws.Range('H116').Value = 1
$ws.Range('I116').Value = 0
$ws.Range('J116').Value = 'H'
$ws.Range('K116').Value = 2
$ws.Range('L116').Value = 3.1
$ws.Range('M116').Value = 3.6
$ws.Range('N116').Value = 2.5
$ws.Range('O116').Value = 3.2
$ws.Range('P116').Value = 2.625
$ws.Range('Q116').Value = 0
$ws.Range('R116').Value = 1.875
$ws.Range('S116').Value = 1.975
$ws.Range('T116').Value = 2.25
$ws.Range('U116').Value = 2.025
$ws.Range('V116').Value = 1.825
$ws.Range('W116').Value = 1.5
$ws.Range('X116').Value = -1
$ws.Range('Y116').Value = -1
$ws.Range('Z116').Value = 0.875
$ws.Range('AA116').Value = -1
$ws.Range('AB116').Value = -1
$ws.Range('AC116').Value = 0.825

# Row 117
$ws.Range('B117').Value = 3373722
$ws.Range('C117').Value = 'Uruguay Primera División'
$ws.Range('D117').Value = 'Uruguay Clausura'
$ws.Range('E117').Value = 44283.83333333334
$ws.Range('F117').Value = 'Deportivo Maldonado'
$ws.Range('G117').Value = 'Nacional De Football'
$ws.Range('H117').Value = 1
$ws.Range('I117').Value = 2
$ws.Range('J117').Value = 'A'
$ws.Range('K117').Value = 3.75
$ws.Range('L117').Value = 3.6
$ws.Range('M117').Value = 1.8
$ws.Range('N117').Value = 5.5
$ws.Range('O117').Value = 4.5
$ws.Range('P117').Value = 1.5
$ws.Range('Q117').Value = 1
$ws.Range('R117').Value = 1.9
$ws.Range('S117').Value = 1.95
$ws.Range('T117').Value = 2.75
$ws.Range('U117').Value = 1.85
$ws.Range('V117').Value = 2
$ws.Range('W117').Value = -1
$ws.Range('X117').Value = -1
$ws.Range('Y117').Value = 0.5
$ws.Range('Z117').Value = 0
$ws.Range('AA117').Value = -0
$ws.Range('AB117').Value = 0.425
$ws.Range('AC117').Value = -0.5

# Row 355
$ws.Range('B355').Value = 4402274
$ws.Range('C355').Value = 'Uruguay Primera División'
$ws.Range('D355').Value = 'Uruguay Clausura'
$ws.Range('E355').Value = 44534.70833333334
$ws.Range('F355').Value = 'Montevideo Wanderers'
$ws.Range('G355').Value = 'Atletico Fenix Montevideo'
$ws.Range('H355').Value = 1
$ws.Range('I355').Value = 1
$ws.Range('J355').Value = 'D'
$ws.Range('K355').Value = 1.909
$ws.Range('L355').Value = 3.25
$ws.Range('M355').Value = 3.6
$ws.Range('N355').Value = 1.909
$ws.Range('O355').Value = 3.3
$ws.Range('P355').Value = 3.6
$ws.Range('Q355').Value = -0.5
$ws.Range('R355').Value = 1.95
$ws.Range('S355').Value = 1.9
$ws.Range('T355').Value = 2.5
$ws.Range('U355').Value = 1.9
$ws.Range('V355').Value = 1.95
$ws.Range('W355').Value = -1
$ws.Range('X355').Value = 2.3
$ws.Range('Y355').Value = -1
$ws.Range('Z355').Value = -1
$ws.Range('AA355').Value = 0.8999999999999999
$ws.Range('AB355').Value = -1
$ws.Range('AC355').Value = 0.95

# Row 356
$ws.Range('B356').Value = 4402276
$ws.Range('C356').Value = 'Uruguay Primera División'
$ws.Range('D356').Value = 'Uruguay Clausura'
$ws.Range('E356').Value = 44534.70833333334
$ws.Range('F356').Value = 'Penarol'
$ws.Range('G356').Value = 'IA Sud America'
$ws.Range('H356').Value = 3
$ws.Range('I356').Value = 1
$ws.Range('J356').Value = 'H'
$ws.Range('K356').Value = 1.333
$ws.Range('L356').Value = 4.5
$ws.Range('M356').Value = 7.5
$ws.Range('N356').Value = 1.25
$ws.Range('O356').Value = 5
$ws.Range('P356').Value = 10
$ws.Range('Q356').Value = -1.5
$ws.Range('R356').Value = 1.85
$ws.Range('S356').Value = 2
$ws.Range('T356').Value = 3
$ws.Range('U356').Value = 1.9
$ws.Range('V356').Value = 1.95
$ws.Range('W356').Value = 0.25
$ws.Range('X356').Value = -1
$ws.Range('Y356').Value = -1
$ws.Range('Z356').Value = 0.8500000000000001
$ws.Range('AA356').Value = -1
$ws.Range('AB356').Value = 0.8999999999999999
$ws.Range('AC356').Value = -1

# Row 357
$ws.Range('B357').Value = 4402277
$ws.Range('C357').Value = 'Uruguay Primera División'
$ws.Range('D357').Value = 'Uruguay Clausura'
$ws.Range('E357').Value = 44534.70833333334
$ws.Range('F357').Value = 'Torque'
$ws.Range('G357').Value = 'Club Atletico Progreso'
$ws.Range('H357').Value = 1
$ws.Range('I357').Value = 0
$ws.Range('J357').Value = 'H'
$ws.Range('K357').Value = 2.05
$ws.Range('L357').Value = 3.4
$ws.Range('M357').Value = 3.1
$ws.Range('N357').Value = 1.909
$ws.Range('O357').Value = 3.5
$ws.Range('P357').Value = 3.3
$ws.Range('Q357').Value = -0.5
$ws.Range('R357').Value = 1.975
$ws.Range('S357').Value = 1.875
$ws.Range('T357').Value = 2.5
$ws.Range('U357').Value = 1.975
$ws.Range('V357').Value = 1.875
$ws.Range('W357').Value = 0.909
$ws.Range('X357').Value = -1
$ws.Range('Y357').Value = -1
$ws.Range('Z357').Value = 0.9750000000000001
$ws.Range('AA357').Value = -1
$ws.Range('AB357').Value = -1
$ws.Range('AC357').Value = 0.875

# Row 358
$ws.Range('B358').Value = 4402275
$ws.Range('C358').Value = 'Uruguay Primera División'
$ws.Range('D358').Value = 'Uruguay Clausura'
$ws.Range('E358').Value = 44534.70833333334
$ws.Range('F358').Value = 'Rentistas'
$ws.Range('G358').Value = 'Plaza Colonia'
$ws.Range('H358').Value = 2
$ws.Range('I358').Value = 3
$ws.Range('J358').Value = 'A'
$ws.Range('K358').Value = 2.9
$ws.Range('L358').Value = 3.2
$ws.Range('M358').Value = 2.25
$ws.Range('N358').Value = 2.6
$ws.Range('O358').Value = 3.2
$ws.Range('P358').Value = 2.5
$ws.Range('Q358').Value = 0
$ws.Range('R358').Value = 1.925
$ws.Range('S358').Value = 1.875
$ws.Range('T358').Value = 2.25
$ws.Range('U358').Value = 2
$ws.Range('V358').Value = 1.8
$ws.Range('W358').Value = -1
$ws.Range('X358').Value = -1
$ws.Range('Y358').Value = 1.5
$ws.Range('Z358').Value = -1
$ws.Range('AA358').Value = 0.875
$ws.Range('AB358').Value = 1
$ws.Range('AC358').Value = -1

# Row 359
$ws.Range('B359').Value = 4402278
$ws.Range('C359').Value = 'Uruguay Primera División'
$ws.Range('D359').Value = 'Uruguay Clausura'
$ws.Range('E359').Value = 44534.70833333334
$ws.Range('F359').Value = 'Nacional De Football'
$ws.Range('G359').Value = 'CA River Plate'
$ws.Range('H359').Value = 4
$ws.Range('I359').Value = 2
$ws.Range('J359').Value = 'H'
$ws.Range('K359').Value = 1.615
$ws.Range('L359').Value = 3.5
$ws.Range('M359').Value = 5
$ws.Range('N359').Value = 1.55
$ws.Range('O359').Value = 3.6
$ws.Range('P359').Value = 5.5
$ws.Range('Q359').Value = -0.75
$ws.Range('R359').Value = 1.8
$ws.Range('S359').Value = 2.05
$ws.Range('T359').Value = 2.75
$ws.Range('U359').Value = 1.975
$ws.Range('V359').Value = 1.875
$ws.Range('W359').Value = 0.55
$ws.Range('X359').Value = -1
$ws.Range('Y359').Value = -1
$ws.Range('Z359').Value = 0.8
$ws.Range('AA359').Value = -1
$ws.Range('AB359').Value = 0.9750000000000001
$ws.Range('AC359').Value = -1

# Row 476
$ws.Range('B476').Value = 4940699
$ws.Range('C476').Value = 'Uruguay Primera División'
$ws.Range('D476').Value = 'Uruguay Apertura'
$ws.Range('E476').Value = 44716.66666666666
$ws.Range('F476').Value = 'Deportivo Maldonado'
$ws.Range('G476').Value = 'Albion FC'
$ws.Range('H476').Value = 0
$ws.Range('I476').Value = 0
$ws.Range('J476').Value = 'D'
$ws.Range('K476').Value = 1.615
$ws.Range('L476').Value = 3.6
$ws.Range('M476').Value = 4.75
$ws.Range('N476').Value = 1.6
$ws.Range('O476').Value = 4
$ws.Range('P476').Value = 4.75
$ws.Range('Q476').Value = -0.75
$ws.Range('R476').Value = 1.825
$ws.Range('S476').Value = 2.025
$ws.Range('T476').Value = 2.75
$ws.Range('U476').Value = 2
$ws.Range('V476').Value = 1.85
$ws.Range('W476').Value = -1
$ws.Range('X476').Value = 3
$ws.Range('Y476').Value = -1
$ws.Range('Z476').Value = -1
$ws.Range('AA476').Value = 1.025
$ws.Range('AB476').Value = -1
$ws.Range('AC476').Value = 0.8500000000000001

# Row 477
$ws.Range('B477').Value = 4939377
$ws.Range('C477').Value = 'Uruguay Primera División'
$ws.Range('D477').Value = 'Uruguay Apertura'
$ws.Range('E477').Value = 44716.66666666666
$ws.Range('F477').Value = 'Penarol'
$ws.Range('G477').Value = 'Cerro Largo'
$ws.Range('H477').Value = 0
$ws.Range('I477').Value = 1
$ws.Range('J477').Value = 'A'
$ws.Range('K477').Value = 1.444
$ws.Range('L477').Value = 3.75
$ws.Range('M477').Value = 7
$ws.Range('N477').Value = 1.533
$ws.Range('O477').Value = 3.4
$ws.Range('P477').Value = 6
$ws.Range('Q477').Value = -1
$ws.Range('R477').Value = 2.05
$ws.Range('S477').Value = 1.8
$ws.Range('T477').Value = 2.25
$ws.Range('U477').Value = 1.925
$ws.Range('V477').Value = 1.925
$ws.Range('W477').Value = -1
$ws.Range('X477').Value = -1
$ws.Range('Y477').Value = 5
$ws.Range('Z477').Value = -1
$ws.Range('AA477').Value = 0.8
$ws.Range('AB477').Value = -1
$ws.Range('AC477').Value = 0.925

# Row 838
$ws.Range('B838').Value = 7013409
$ws.Range('C838').Value = 'Uruguay Primera División'
$ws.Range('D838').Value = 'Uruguay Clausura'
$ws.Range('E838').Value = 45267.70833333334
$ws.Range('F838').Value = 'Nacional De Football'
$ws.Range('G838').Value = 'Torque'
$ws.Range('H838').Value = 1
$ws.Range('I838').Value = 1
$ws.Range('J838').Value = 'D'
$ws.Range('K838').Value = 1.666
$ws.Range('L838').Value = 3.9
$ws.Range('M838').Value = 4.5
$ws.Range('N838').Value = 1.615
$ws.Range('O838').Value = 4
$ws.Range('P838').Value = 4.75
$ws.Range('Q838').Value = -0.75
$ws.Range('R838').Value = 1.8
$ws.Range('S838').Value = 2.05
$ws.Range('T838').Value = 2.75
$ws.Range('U838').Value = 1.95
$ws.Range('V838').Value = 1.9
$ws.Range('W838').Value = -1
$ws.Range('X838').Value = 3
$ws.Range('Y838').Value = -1
$ws.Range('Z838').Value = -1
$ws.Range('AA838').Value = 1.05
$ws.Range('AB838').Value = -1
$ws.Range('AC838').Value = 0.8999999999999999

# Row 839
$ws.Range('B839').Value = 7013885
$ws.Range('C839').Value = 'Uruguay Primera División'
$ws.Range('D839').Value = 'Uruguay Clausura'
$ws.Range('E839').Value = 45267.70833333334
$ws.Range('F839').Value = 'La Luz'
$ws.Range('G839').Value = 'Atletico Fenix Montevideo'
$ws.Range('H839').Value = 0
$ws.Range('I839').Value = 2
$ws.Range('J839').Value = 'A'
$ws.Range('K839').Value = 3
$ws.Range('L839').Value = 3
$ws.Range('M839').Value = 2.4
$ws.Range('N839').Value = 2.9
$ws.Range('O839').Value = 2.75
$ws.Range('P839').Value = 2.6
$ws.Range('Q839').Value = 0
$ws.Range('R839').Value = 2.025
$ws.Range('S839').Value = 1.825
$ws.Range('T839').Value = 2
$ws.Range('U839').Value = 2.025
$ws.Range('V839').Value = 1.825
$ws.Range('W839').Value = -1
$ws.Range('X839').Value = -1
$ws.Range('Y839').Value = 1.6
$ws.Range('Z839').Value = -1
$ws.Range('AA839').Value = 0.825
$ws.Range('AB839').Value = 0
$ws.Range('AC839').Value = -0
